$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 5
$ws.Range("G5").Value = 1.25
$ws.Range("H5").Value = 5.75
$ws.Range("I5").Value = 9.5
$ws.Range("K5").Value = 34
$ws.Range("R5").Value = 1.53
$ws.Range("S5").Value = 2.38
$ws.Range("T5").Value = 17
$ws.Range("U5").Value = 11
$ws.Range("V5").Value = 11
$ws.Range("W5").Value = 11
$ws.Range("Y5").Value = 19
$ws.Range("Z5").Value = 34
$ws.Range("AA5").Value = 13
$ws.Range("AB5").Value = 17
$ws.Range("AE5").Value = 51
$ws.Range("AF5").Value = 29
$ws.Range("AG5").Value = 101
$ws.Range("AH5").Value = 51
$ws.Range("AI5").Value = 41
$ws.Range("AJ5").Value = 101

# Row 6
$ws.Range("G6").Value = 2.18
$ws.Range("L6").Value = 1.29
$ws.Range("M6").Value = 3.4
$ws.Range("N6").Value = 1.87
$ws.Range("O6").Value = 1.87
$ws.Range("P6").Value = 1.4
$ws.Range("Q6").Value = 2.82
$ws.Range("S6").Value = 2.02
$ws.Range("U6").Value = 11.5
$ws.Range("Y6").Value = 29
$ws.Range("AA6").Value = 6.9
$ws.Range("AB6").Value = 14.5
$ws.Range("AC6").Value = 65
$ws.Range("AD6").Value = 9.75
$ws.Range("AE6").Value = 17.5
$ws.Range("AH6").Value = 29
$ws.Range("AI6").Value = 37
$ws.Range("AJ6").Value = 500

# Row 13
$ws.Range("J13").Value = 1.03
$ws.Range("K13").Value = 15
$ws.Range("N13").Value = 1.7
$ws.Range("O13").Value = 2.1

# Row 14
$ws.Range("N14").Value = 2.2
$ws.Range("O14").Value = 1.65

# Row 20
$ws.Range("H20").Value = 3.2
$ws.Range("I20").Value = 3.5
$ws.Range("N20").Value = 2.3
$ws.Range("O20").Value = 1.6
$ws.Range("R20").Value = 1.91
$ws.Range("S20").Value = 1.8
$ws.Range("T20").Value = 6.5
$ws.Range("U20").Value = 9.5
$ws.Range("W20").Value = 21
$ws.Range("Y20").Value = 34
$ws.Range("Z20").Value = 8
$ws.Range("AB20").Value = 17
$ws.Range("AD20").Value = 9
$ws.Range("AF20").Value = 13
$ws.Range("AG20").Value = 41
$ws.Range("AJ20").Value = 351

# Row 25
$ws.Range("G25").Value = 1.42
$ws.Range("L25").Value = 1.25
$ws.Range("M25").Value = 3.75

# Row 27
$ws.Range("G27").Value = 2.3
$ws.Range("H27").Value = 2.8
$ws.Range("N27").Value = 2.77
$ws.Range("O27").Value = 1.33
$ws.Range("Q27").Value = 2.05
$ws.Range("R27").Value = 2.27
$ws.Range("V27").Value = 10.5
$ws.Range("W27").Value = 23
$ws.Range("X27").Value = 26
$ws.Range("Y27").Value = 55
$ws.Range("AA27").Value = 5.9
$ws.Range("AB27").Value = 23
$ws.Range("AD27").Value = 6.6
$ws.Range("AF27").Value = 14
$ws.Range("AH27").Value = 45
$ws.Range("AI27").Value = 80

# Row 29
$ws.Range("J29").Value = 1.02
$ws.Range("K29").Value = 19

# Row 31
$ws.Range("G31").Value = 3.75
$ws.Range("I31").Value = 2.05
$ws.Range("J31").Value = 1.05
$ws.Range("K31").Value = 11
$ws.Range("R31").Value = 1.75
$ws.Range("S31").Value = 2
$ws.Range("T31").Value = 12
$ws.Range("Z31").Value = 11

# Row 32
$ws.Range("G32").Value = 1.65
$ws.Range("H32").Value = 3.6
$ws.Range("I32").Value = 4.6
$ws.Range("P32").Value = 1.44
$ws.Range("Q32").Value = 2.62
$ws.Range("R32").Value = 1.93
$ws.Range("S32").Value = 1.78
$ws.Range("T32").Value = 6.3
$ws.Range("U32").Value = 7.3
$ws.Range("W32").Value = 12
$ws.Range("AA32").Value = 7.1
$ws.Range("AB32").Value = 18
$ws.Range("AC32").Value = 90
$ws.Range("AD32").Value = 11.75
$ws.Range("AF32").Value = 15.5
$ws.Range("AI32").Value = 55
$ws.Range("AJ32").Value = 800
